$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 13159274
$ws.Range("I80").Value = 22728322
$ws.Range("J80").Value = 1833.375
$ws.Range("K80").Value = 68184966
$ws.Range("L80").Value = 5500.125
$ws.Range("M80").Value = -68183968
$ws.Range("N80").Value = -7496.125
$ws.Range("H83").Value = 13159274
$ws.Range("I83").Value = 22728322
$ws.Range("J83").Value = 1833.375
$ws.Range("K83").Value = 204554898
$ws.Range("L83").Value = 16500.375
$ws.Range("M83").Value = -204549906
$ws.Range("N83").Value = -26484.375
$ws.Range("H86").Value = 47620984
$ws.Range("I86").Value = 52633590
$ws.Range("J86").Value = 1222
$ws.Range("K86").Value = 52633590
$ws.Range("L86").Value = 1222
$ws.Range("M86").Value = -52632467
$ws.Range("N86").Value = -3468
$ws.Range("H88").Value = 4343.615
$ws.Range("I88").Value = 671.625
$ws.Range("J88").Value = 10218.8
$ws.Range("K88").Value = 671.625
$ws.Range("L88").Value = 10218.8
$ws.Range("M88").Value = -265.625
$ws.Range("N88").Value = -11030.8
$ws.Range("H89").Value = 47620984
$ws.Range("I89").Value = 52633590
$ws.Range("J89").Value = 1222
$ws.Range("K89").Value = 263167950
$ws.Range("L89").Value = 6110
$ws.Range("M89").Value = -263162334
$ws.Range("N89").Value = -17342
$ws.Range("H91").Value = 4343.615
$ws.Range("I91").Value = 671.625
$ws.Range("J91").Value = 10218.8
$ws.Range("K91").Value = 671.625
$ws.Range("L91").Value = 10218.8
$ws.Range("M91").Value = 732.375
$ws.Range("N91").Value = -13026.8
$ws.Range("H92").Value = 34483140
$ws.Range("I92").Value = 35714670
$ws.Range("J92").Value = 300
$ws.Range("K92").Value = 35714670
$ws.Range("L92").Value = 300
$ws.Range("M92").Value = -35713422
$ws.Range("N92").Value = -2796
$ws.Range("H99").Value = 6035.385
$ws.Range("J99").Value = 12596.4
$ws.Range("L99").Value = 37789.2
$ws.Range("N99").Value = -40785.2
$ws.Range("H103").Value = 779.0909
$ws.Range("J103").Value = 883.875
$ws.Range("L103").Value = 2651.625
$ws.Range("N103").Value = -3823.625
$ws.Range("H132").Value = 191446.53
$ws.Range("I132").Value = 208451.11
$ws.Range("K132").Value = 625353.33
$ws.Range("M132").Value = -622823.33
$ws.Range("H138").Value = 4832.235
$ws.Range("J138").Value = 5270.012
$ws.Range("L138").Value = 15810.036
$ws.Range("N138").Value = -26090.036

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 21047.24
$ws.Range("I2").Value = 27085.264
$ws.Range("K2").Value = 27085.264
$ws.Range("M2").Value = -26972.264
$ws.Range("H32").Value = 15714.233
$ws.Range("I32").Value = 15078.691
$ws.Range("J32").Value = 21434.111
$ws.Range("K32").Value = 15078.691
$ws.Range("L32").Value = 21434.111
$ws.Range("M32").Value = -14791.691
$ws.Range("N32").Value = -22008.111
$ws.Range("H45").Value = 26892.203
$ws.Range("I45").Value = 31921.426
$ws.Range("K45").Value = 31921.426
$ws.Range("M45").Value = -31544.426
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("N55").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("N96").Value = 0
$ws.Range("H110").Value = 1438.7222
$ws.Range("I110").Value = 1435.1177
$ws.Range("K110").Value = 1435.1177
$ws.Range("M110").Value = 609.8823
$ws.Range("H116").Value = 21047.24
$ws.Range("I116").Value = 27085.264
$ws.Range("K116").Value = 27085.264
$ws.Range("M116").Value = -24791.264
$ws.Range("H122").Value = 4438.5835
$ws.Range("I122").Value = 4237.273
$ws.Range("K122").Value = 12711.819
$ws.Range("M122").Value = -10261.819
$ws.Range("H132").Value = 10198.453
$ws.Range("I132").Value = 10389.75
$ws.Range("K132").Value = 31169.25
$ws.Range("M132").Value = -28639.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 21047.24
$ws.Range("I3").Value = 27085.264
$ws.Range("K3").Value = 27085.264
$ws.Range("M3").Value = -26971.264
$ws.Range("H134").Value = 5149.9375
$ws.Range("I134").Value = 3034.6428
$ws.Range("K134").Value = 9103.928400000001
$ws.Range("M134").Value = -6568.928400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1415.25
$ws.Range("I16").Value = 1020.3333
$ws.Range("J16").Value = 2600
$ws.Range("K16").Value = 1020.3333
$ws.Range("L16").Value = 2600
$ws.Range("M16").Value = -733.3333
$ws.Range("N16").Value = -3174
$ws.Range("H22").Value = 504.1
$ws.Range("I22").Value = 372.5
$ws.Range("J22").Value = 591.8333
$ws.Range("K22").Value = 372.5
$ws.Range("L22").Value = 591.8333
$ws.Range("M22").Value = -22.5
$ws.Range("N22").Value = -1291.8333
$ws.Range("H58").Value = 2123.7
$ws.Range("I58").Value = 1981.8889
$ws.Range("K58").Value = 1981.8889
$ws.Range("M58").Value = -1778.8889
$ws.Range("H99").Value = 7034.25
$ws.Range("I99").Value = 3724.75
$ws.Range("J99").Value = 8137.4165
$ws.Range("K99").Value = 3724.75
$ws.Range("L99").Value = 8137.4165
$ws.Range("M99").Value = -2226.75
$ws.Range("N99").Value = -11133.4165
$ws.Range("H105").Value = 1598.5
$ws.Range("I105").Value = 1446
$ws.Range("K105").Value = 1446
$ws.Range("M105").Value = 301
$ws.Range("H107").Value = 940.95654
$ws.Range("I107").Value = 784.58826
$ws.Range("K107").Value = 784.58826
$ws.Range("M107").Value = 1135.41174
$ws.Range("H113").Value = 1415.25
$ws.Range("I113").Value = 1020.3333
$ws.Range("J113").Value = 2600
$ws.Range("K113").Value = 1020.3333
$ws.Range("L113").Value = 2600
$ws.Range("M113").Value = 1149.6667
$ws.Range("N113").Value = -6940
$ws.Range("H126").Value = 7034.25
$ws.Range("I126").Value = 3724.75
$ws.Range("J126").Value = 8137.4165
$ws.Range("K126").Value = 11174.25
$ws.Range("L126").Value = 24412.2495
$ws.Range("M126").Value = -8704.25
$ws.Range("N126").Value = -29352.2495
$ws.Range("H136").Value = 2123.7
$ws.Range("I136").Value = 1981.8889
$ws.Range("K136").Value = 5945.6667
$ws.Range("M136").Value = -3395.6667
$ws.Range("H138").Value = 97583.92
$ws.Range("J138").Value = 97583.92
$ws.Range("L138").Value = 97583.92
$ws.Range("N138").Value = -107863.92

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 678.35
$ws.Range("I12").Value = 1461.625
$ws.Range("K12").Value = 4384.875
$ws.Range("M12").Value = -4211.875
$ws.Range("H61").Value = 526.55554
$ws.Range("J61").Value = 1019.3333
$ws.Range("L61").Value = 3057.9999
$ws.Range("N61").Value = -3487.9999
$ws.Range("H134").Value = 6035.2
$ws.Range("I134").Value = 5594.6665
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 16783.9995
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -11713.9995
$ws.Range("N134").Value = -40140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 267.6
$ws.Range("I2").Value = 127.875
$ws.Range("J2").Value = 360.75
$ws.Range("K2").Value = 127.875
$ws.Range("L2").Value = 360.75
$ws.Range("M2").Value = -14.875
$ws.Range("N2").Value = -586.75
$ws.Range("H57").Value = 21998.6
$ws.Range("J57").Value = 21998.6
$ws.Range("L57").Value = 21998.6
$ws.Range("N57").Value = -23638.6
$ws.Range("H80").Value = 69699.836
$ws.Range("I80").Value = 97775.91
$ws.Range("K80").Value = 97775.91
$ws.Range("M80").Value = -96777.91
$ws.Range("H83").Value = 69699.836
$ws.Range("I83").Value = 97775.91
$ws.Range("K83").Value = 488879.55
$ws.Range("M83").Value = -483887.55
$ws.Range("H132").Value = 6142.241
$ws.Range("J132").Value = 5721.75
$ws.Range("L132").Value = 17165.25
$ws.Range("N132").Value = -22225.25
$ws.Range("H138").Value = 91388.5
$ws.Range("J138").Value = 91388.5
$ws.Range("L138").Value = 91388.5
$ws.Range("N138").Value = -101668.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4577.143
$ws.Range("I7").Value = 3331.2222
$ws.Range("K7").Value = 3331.2222
$ws.Range("M7").Value = -3219.2222
$ws.Range("H126").Value = 4577.143
$ws.Range("I126").Value = 3331.2222
$ws.Range("K126").Value = 9993.6666
$ws.Range("M126").Value = -7523.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 17448.467
$ws.Range("J81").Value = 20094
$ws.Range("L81").Value = 40188
$ws.Range("N81").Value = -42310
$ws.Range("H84").Value = 17448.467
$ws.Range("J84").Value = 20094
$ws.Range("L84").Value = 200940
$ws.Range("N84").Value = -211548
